$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q3" sheet right after "总计" (i.e. right
#    before "2022-Q2"), by copying the "2022-Q2" sheet so that it
#    inherits the exact same formatting/styles as its siblings.
# ------------------------------------------------------------------
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($srcQ2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Overwrite the copied values with the real 2022-Q3 figures.
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'513030"
$q3.Range("C2").Value = "华安国际龙头（DAX）ETF（QDII）"
$q3.Range("D2").Value = "'5.54"
$q3.Range("E2").Value = "'93.57"
$q3.Range("F2").Value = "'5.59"
$q3.Range("G2").Value = "'0.3097"
$q3.Range("H2").Value = 5

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'006282"
$q3.Range("C3").Value = "上投摩根欧洲动力策略股票（QDII）"
$q3.Range("D3").Value = "'0.41"
$q3.Range("E3").Value = "'91.47"
$q3.Range("F3").Value = "'1.84"
$q3.Range("G3").Value = "'0.0075"
$q3.Range("H3").Value = 9

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: a new row for 2022-Q3 is
#    inserted at row 2, pushing the previously existing rows down by
#    one (their index column A is recomputed as a 0-based sequence).
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.32

$total.Range("A3").Value = 1
$total.Range("D3").Value = 0.35

$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7

Write-Host "done"
